{"js": "const pairs = [\n  [\"2023-08-13 Sunday\", \"2023-08-14 Monday\"],\n  [\"19-14=\", \"9+64=\"],\n  [\"47+40=\", \"44-14=\"],\n  [\"14+69=\", \"5+48=\"],\n  [\"48+49=\", \"68-59=\"],\n  [\"80-19=\", \"91-3=\"],\n  [\"70+23=\", \"35-27=\"],\n  [\"38-19=\", \"82+7=\"],\n  [\"74-45=\", \"38+27=\"],\n  [\"60-14=\", \"18-6=\"],\n  [\"46-39=\", \"90+6=\"],\n  [\"14+9=\", \"87-21=\"],\n  [\"38+8=\", \"43+51=\"],\n  [\"51-45=\", \"71+23=\"],\n  [\"22+55=\", \"74-22=\"],\n  [\"93-49=\", \"97-8=\"],\n  [\"49+23=\", \"74-30=\"],\n  [\"42+39=\", \"90-67=\"],\n  [\"49+10=\", \"47-13=\"],\n  [\"7+15=\", \"31+53=\"],\n  [\"97-67=\", \"32-14=\"],\n  [\"44-10=\", \"95-55=\"],\n  [\"43-9=\", \"38+22=\"],\n  [\"83-45=\", \"75-48=\"],\n  [\"70-11=\", \"3+9=\"],\n  [\"43-6=\", \"13+29=\"],\n  [\"35+35=\", \"17+50=\"],\n  [\"26+20=\", \"46+35=\"],\n  [\"66-14=\", \"23-12=\"],\n  [\"6+45=\", \"8+22=\"],\n  [\"40+47=\", \"11+16=\"],\n  [\"56-21=\", \"79-48=\"],\n  [\"77-7=\", \"72-56=\"],\n  [\"47+36=\", \"59+15=\"],\n  [\"51-30=\", \"50-1=\"],\n  [\"29+65=\", \"43-22=\"],\n  [\"10+48=\", \"26+14=\"],\n  [\"4+13=\", \"18+71=\"],\n  [\"8+73=\", \"64+8=\"],\n  [\"6+33=\", \"9+9=\"],\n  [\"4+63=\", \"55+30=\"],\n  [\"18+67=\", \"14+59=\"],\n  [\"19+41=\", \"94-10=\"],\n  [\"64-10=\", \"88-13=\"],\n  [\"89-8=\", \"82-60=\"],\n  [\"0+29=\", \"2+22=\"],\n  [\"43+12=\", \"58+30=\"],\n  [\"80-70=\", \"74-49=\"],\n  [\"25+6=\", \"93+2=\"],\n  [\"73-73=\", \"91-1=\"],\n  [\"40-2=\", \"12+0=\"],\n  [\"61+19=\", \"25+31=\"],\n  [\"47+22=\", \"79-4=\"],\n  [\"87-3=\", \"69-51=\"],\n  [\"96-74=\", \"16+66=\"],\n  [\"60-17=\", \"92+1=\"],\n  [\"12+48=\", \"69+4=\"],\n  [\"6+11=\", \"69-11=\"],\n  [\"16+4=\", \"93-88=\"],\n  [\"70+10=\", \"64+14=\"],\n  [\"10+33=\", \"50+8=\"],\n  [\"0+5=\", \"98-59=\"],\n  [\"28+62=\", \"60+34=\"],\n  [\"36+56=\", \"17+40=\"],\n  [\"66-22=\", \"69-59=\"],\n  [\"40+10=\", \"56-10=\"],\n  [\"88-5=\", \"33+58=\"],\n  [\"51-7=\", \"95-75=\"],\n  [\"91-66=\", \"57+13=\"],\n  [\"61+27=\", \"30+31=\"],\n  [\"58+34=\", \"7+45=\"],\n  [\"13+4=\", \"52-18=\"],\n  [\"68+19=\", \"97-41=\"],\n  [\"68-18=\", \"56-7=\"],\n  [\"11-4=\", \"82+4=\"],\n  [\"83-62=\", \"96-53=\"],\n  [\"67-56=\", \"92-31=\"],\n  [\"14+66=\", \"81-25=\"],\n  [\"37-22=\", \"13+8=\"],\n  [\"10+61=\", \"99-77=\"],\n  [\"70-58=\", \"84-44=\"],\n  [\"5+1=\", \"81-10=\"],\n  [\"59-41=\", \"83-13=\"],\n  [\"65-47=\", \"93-10=\"],\n  [\"70-47=\", \"56+23=\"],\n  [\"0+24=\", \"15+10=\"],\n  [\"50-3=\", \"94-5=\"],\n  [\"68-5=\", \"27+67=\"],\n  [\"5+40=\", \"16+42=\"],\n  [\"18+28=\", \"68-49=\"],\n  [\"67-0=\", \"72-25=\"],\n  [\"2+85=\", \"7-4=\"],\n  [\"99-24=\", \"12+71=\"],\n  [\"88-84=\", \"76-6=\"],\n  [\"41+7=\", \"96-95=\"],\n  [\"19+12=\", \"41-14=\"],\n  [\"62-10=\", \"85-7=\"],\n  [\"89-18=\", \"45+22=\"],\n  [\"19+39=\", \"59+9=\"],\n  [\"11+75=\", \"50-8=\"],\n  [\"45-10=\", \"26-1=\"],\n];\n\nconst body = context.document.body;\n\n// Each old string below occurs exactly once in the document, so a literal,\n// case-sensitive search-and-replace (applied one pair at a time) reproduces\n// the diff exactly regardless of processing order.\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@('2023-08-13 Sunday', '2023-08-14 Monday')\n    ,@('19-14=', '9+64=')\n    ,@('47+40=', '44-14=')\n    ,@('14+69=', '5+48=')\n    ,@('48+49=', '68-59=')\n    ,@('80-19=', '91-3=')\n    ,@('70+23=', '35-27=')\n    ,@('38-19=', '82+7=')\n    ,@('74-45=', '38+27=')\n    ,@('60-14=', '18-6=')\n    ,@('46-39=', '90+6=')\n    ,@('14+9=', '87-21=')\n    ,@('38+8=', '43+51=')\n    ,@('51-45=', '71+23=')\n    ,@('22+55=', '74-22=')\n    ,@('93-49=', '97-8=')\n    ,@('49+23=', '74-30=')\n    ,@('42+39=', '90-67=')\n    ,@('49+10=', '47-13=')\n    ,@('7+15=', '31+53=')\n    ,@('97-67=', '32-14=')\n    ,@('44-10=', '95-55=')\n    ,@('43-9=', '38+22=')\n    ,@('83-45=', '75-48=')\n    ,@('70-11=', '3+9=')\n    ,@('43-6=', '13+29=')\n    ,@('35+35=', '17+50=')\n    ,@('26+20=', '46+35=')\n    ,@('66-14=', '23-12=')\n    ,@('6+45=', '8+22=')\n    ,@('40+47=', '11+16=')\n    ,@('56-21=', '79-48=')\n    ,@('77-7=', '72-56=')\n    ,@('47+36=', '59+15=')\n    ,@('51-30=', '50-1=')\n    ,@('29+65=', '43-22=')\n    ,@('10+48=', '26+14=')\n    ,@('4+13=', '18+71=')\n    ,@('8+73=', '64+8=')\n    ,@('6+33=', '9+9=')\n    ,@('4+63=', '55+30=')\n    ,@('18+67=', '14+59=')\n    ,@('19+41=', '94-10=')\n    ,@('64-10=', '88-13=')\n    ,@('89-8=', '82-60=')\n    ,@('0+29=', '2+22=')\n    ,@('43+12=', '58+30=')\n    ,@('80-70=', '74-49=')\n    ,@('25+6=', '93+2=')\n    ,@('73-73=', '91-1=')\n    ,@('40-2=', '12+0=')\n    ,@('61+19=', '25+31=')\n    ,@('47+22=', '79-4=')\n    ,@('87-3=', '69-51=')\n    ,@('96-74=', '16+66=')\n    ,@('60-17=', '92+1=')\n    ,@('12+48=', '69+4=')\n    ,@('6+11=', '69-11=')\n    ,@('16+4=', '93-88=')\n    ,@('70+10=', '64+14=')\n    ,@('10+33=', '50+8=')\n    ,@('0+5=', '98-59=')\n    ,@('28+62=', '60+34=')\n    ,@('36+56=', '17+40=')\n    ,@('66-22=', '69-59=')\n    ,@('40+10=', '56-10=')\n    ,@('88-5=', '33+58=')\n    ,@('51-7=', '95-75=')\n    ,@('91-66=', '57+13=')\n    ,@('61+27=', '30+31=')\n    ,@('58+34=', '7+45=')\n    ,@('13+4=', '52-18=')\n    ,@('68+19=', '97-41=')\n    ,@('68-18=', '56-7=')\n    ,@('11-4=', '82+4=')\n    ,@('83-62=', '96-53=')\n    ,@('67-56=', '92-31=')\n    ,@('14+66=', '81-25=')\n    ,@('37-22=', '13+8=')\n    ,@('10+61=', '99-77=')\n    ,@('70-58=', '84-44=')\n    ,@('5+1=', '81-10=')\n    ,@('59-41=', '83-13=')\n    ,@('65-47=', '93-10=')\n    ,@('70-47=', '56+23=')\n    ,@('0+24=', '15+10=')\n    ,@('50-3=', '94-5=')\n    ,@('68-5=', '27+67=')\n    ,@('5+40=', '16+42=')\n    ,@('18+28=', '68-49=')\n    ,@('67-0=', '72-25=')\n    ,@('2+85=', '7-4=')\n    ,@('99-24=', '12+71=')\n    ,@('88-84=', '76-6=')\n    ,@('41+7=', '96-95=')\n    ,@('19+12=', '41-14=')\n    ,@('62-10=', '85-7=')\n    ,@('89-18=', '45+22=')\n    ,@('19+39=', '59+9=')\n    ,@('11+75=', '50-8=')\n    ,@('45-10=', '26-1=')\n)\n\n# Every \"old\" string below occurs exactly once in the document, so a\n# literal, case-sensitive find-and-replace-all (run once per pair) reproduces\n# the diff exactly, regardless of the order the pairs are processed in.\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    # FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    # MatchAllWordForms, Forward, Wrap(wdFindContinue=1), Format, ReplaceWith,\n    # Replace(wdReplaceAll=2)\n    $result = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $result) {\n        throw \"No match found for: $oldText\"\n    }\n}"}
